# Update the "scraped_at" timestamps (column K) on the "snapshot" sheet
# to reflect the latest scrape run, per the commit:
# "chore(runtime): publish files + archive (2025-10-30 12:27:33)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$newTimestamps = @{
    2 = "2025-10-30T04:26:24.051260+00:00"
    3 = "2025-10-30T04:26:26.082601+00:00"
    4 = "2025-10-30T04:26:26.082617+00:00"
    5 = "2025-10-30T04:26:26.082625+00:00"
    6 = "2025-10-30T04:26:26.082633+00:00"
    7 = "2025-10-30T04:26:26.082640+00:00"
    8 = "2025-10-30T04:26:28.007987+00:00"
    9 = "2025-10-30T04:26:28.008019+00:00"
    10 = "2025-10-30T04:26:30.035558+00:00"
    11 = "2025-10-30T04:26:32.770902+00:00"
    12 = "2025-10-30T04:26:32.770931+00:00"
    13 = "2025-10-30T04:26:35.098490+00:00"
    14 = "2025-10-30T04:26:35.098518+00:00"
    15 = "2025-10-30T04:26:35.098536+00:00"
    16 = "2025-10-30T04:26:35.098552+00:00"
    17 = "2025-10-30T04:26:43.201049+00:00"
    18 = "2025-10-30T04:26:45.956182+00:00"
    19 = "2025-10-30T04:26:48.261303+00:00"
    20 = "2025-10-30T04:26:50.985374+00:00"
    21 = "2025-10-30T04:26:50.985406+00:00"
    22 = "2025-10-30T04:26:50.985424+00:00"
    23 = "2025-10-30T04:26:53.331844+00:00"
    24 = "2025-10-30T04:26:53.331880+00:00"
    25 = "2025-10-30T04:26:53.331900+00:00"
    26 = "2025-10-30T04:26:53.331918+00:00"
    27 = "2025-10-30T04:26:53.331937+00:00"
    28 = "2025-10-30T04:27:00.303160+00:00"
    29 = "2025-10-30T04:27:00.303190+00:00"
    30 = "2025-10-30T04:27:00.303209+00:00"
    31 = "2025-10-30T04:27:00.303227+00:00"
    32 = "2025-10-30T04:27:03.114542+00:00"
    33 = "2025-10-30T04:27:03.114573+00:00"
    34 = "2025-10-30T04:27:03.114592+00:00"
    35 = "2025-10-30T04:27:05.491739+00:00"
    36 = "2025-10-30T04:27:05.491769+00:00"
    37 = "2025-10-30T04:27:05.491788+00:00"
    38 = "2025-10-30T04:27:05.491805+00:00"
    39 = "2025-10-30T04:27:05.491822+00:00"
    40 = "2025-10-30T04:27:05.491838+00:00"
    41 = "2025-10-30T04:27:05.491855+00:00"
    42 = "2025-10-30T04:27:05.491875+00:00"
    43 = "2025-10-30T04:27:05.491891+00:00"
    44 = "2025-10-30T04:27:07.764916+00:00"
    45 = "2025-10-30T04:27:07.764946+00:00"
    46 = "2025-10-30T04:27:12.799825+00:00"
    47 = "2025-10-30T04:27:15.472434+00:00"
    48 = "2025-10-30T04:27:15.472462+00:00"
    49 = "2025-10-30T04:27:15.472479+00:00"
    50 = "2025-10-30T04:27:15.472494+00:00"
}

foreach ($row in $newTimestamps.Keys) {
    $ws.Cells.Item([int]$row, 11).Value = $newTimestamps[$row]
}
